$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark newly-absent sessions with "A" for the affected students/sessions.
# This also drives the "Total Absence" COUNTIF formulas in column E to recompute.
$ws.Range("V7").Value = "A"
$ws.Range("Y7").Value = "A"
$ws.Range("X8").Value = "A"
$ws.Range("X9").Value = "A"
$ws.Range("Y11").Value = "A"
$ws.Range("V12").Value = "A"
$ws.Range("X12").Value = "A"
$ws.Range("Y12").Value = "A"
$ws.Range("V13").Value = "A"
$ws.Range("Y13").Value = "A"
$ws.Range("V15").Value = "A"
$ws.Range("W15").Value = "A"
$ws.Range("X16").Value = "A"
$ws.Range("X19").Value = "A"
$ws.Range("Y19").Value = "A"
$ws.Range("V20").Value = "A"
$ws.Range("V21").Value = "A"
$ws.Range("W22").Value = "A"
$ws.Range("X26").Value = "A"
$ws.Range("Y26").Value = "A"
$ws.Range("W29").Value = "A"
$ws.Range("Y29").Value = "A"
$ws.Range("Y32").Value = "A"
$ws.Range("X33").Value = "A"
$ws.Range("X35").Value = "A"
$ws.Range("Y35").Value = "A"
$ws.Range("W37").Value = "A"
$ws.Range("W38").Value = "A"
$ws.Range("X38").Value = "A"
$ws.Range("W40").Value = "A"
$ws.Range("X40").Value = "A"
$ws.Range("Y40").Value = "A"
$ws.Range("V41").Value = "A"
$ws.Range("W41").Value = "A"
$ws.Range("W42").Value = "A"
$ws.Range("Y42").Value = "A"
$ws.Range("X43").Value = "A"
$ws.Range("Y43").Value = "A"
$ws.Range("X45").Value = "A"
$ws.Range("X46").Value = "A"
$ws.Range("V47").Value = "A"
$ws.Range("V48").Value = "A"
$ws.Range("Y50").Value = "A"
$ws.Range("V51").Value = "A"
$ws.Range("V52").Value = "A"
$ws.Range("V54").Value = "A"
$ws.Range("X55").Value = "A"
$ws.Range("Y55").Value = "A"
$ws.Range("Y56").Value = "A"
$ws.Range("X59").Value = "A"
$ws.Range("Y59").Value = "A"
$ws.Range("V60").Value = "A"
$ws.Range("W60").Value = "A"
$ws.Range("Y60").Value = "A"
$ws.Range("Y61").Value = "A"
$ws.Range("X64").Value = "A"
$ws.Range("Y65").Value = "A"
